# Module 3 - Git and GitHub : "New and Improved Course" commit
#
# Adds extra guidance to the "git push origin master" slide (Workflow /
# Step 4 - Push changes to central repository) explaining how to inspect
# a remote with `git remote show origin` and how to rename a remote with
# `git rename <old> <new>`. The new bullets land right after the
# "origin is a remote pointing to the URL of the repository." paragraph
# and before the "master is the branch ..." paragraph.

$p = $ppt.ActivePresentation

# Find the slide that contains the "git push origin master" content
# placeholder, rather than hard-coding a slide index.
$targetSlide = $null
$targetShape = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    for ($j = 1; $j -le $candidate.Shapes.Count; $j++) {
        $shp = $candidate.Shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -like "*git push origin master*") {
                $targetSlide = $candidate
                $targetShape = $shp
            }
        }
    }
}

$shape = $targetShape
$tf = $shape.TextFrame
$tr = $tf.TextRange

# New bullets are being inserted, so let PowerPoint shrink the text to
# fit instead of overflowing the shape.
$tf.AutoSize = 2

$enDash = [char]0x2013

# Paragraph 3 is "origin is a remote pointing to the URL of the repository."
# The new bullets must land right after it (and before "master is the
# branch ...", currently paragraph 4).
$anchor = $tr.Paragraphs(3, 1)

# --- "You can use git remote show origin ... (Or git remote -v.)" ---
[void]$anchor.InsertAfter("`rYou can use git remote show origin to see details on that. (Or git remote " + $enDash + "v.)")
$para = $tr.Paragraphs(4, 1)
$para.IndentLevel = 2

$run = $para.Characters(13, 23)          # "git remote show origin "
$run.Font.Name = "Consolas"
$run.Font.Size = 18

$run = $para.Characters(64, 13)          # "git remote -v"
$run.Font.Name = "Consolas"
$run.Font.Size = 18

# --- "You can rename your remote using git rename <old> <new>. For example:" ---
[void]$para.InsertAfter("`rYou can rename your remote using git rename <old> <new>. For example:")
$para = $tr.Paragraphs(5, 1)
$para.IndentLevel = 2

$run = $para.Characters(34, 22)          # "git rename <old> <new>"
$run.Font.Name = "Consolas"
$run.Font.Size = 18

# --- "<tab>git rename origin ghrepo" example command, no bullet ---
[void]$para.InsertAfter("`r`tgit rename origin ghrepo")
$para = $tr.Paragraphs(6, 1)
$para.IndentLevel = 2
$para.ParagraphFormat.Bullet.Type = 0    # ppBulletNone -> <a:buNone/>

$run = $para.Characters(2, 18)           # "git rename origin "
$run.Font.Name = "Consolas"
$run.Font.Size = 18

$run = $para.Characters(20, 6)           # "ghrepo"
$run.Font.Name = "Consolas"
$run.Font.Size = 18
